$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 36: dmqm_seminar entry
$ws.Range("D36").Value = "Unifying contrastive learning and clustering"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/386"

# Row 50: etc entry
$ws.Range("D50").Value = "NEOM city"
$ws.Range("E50").Value = "http://incredible.egloos.com/7573911"

# Row 51: bskyvsion entry
$ws.Range("D51").Value = "해피해킹 키보드 아이폰 블루투스 연결하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%ED%95%B4%ED%94%BC%ED%95%B4%ED%82%B9-%ED%82%A4%EB%B3%B4%EB%93%9C-%EC%95%84%EC%9D%B4%ED%8F%B0-%EB%B8%94%EB%A3%A8%ED%88%AC%EC%8A%A4-%EC%97%B0%EA%B2%B0%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95"
